$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 4: B4, C4, D4)
$ws.Range("B4").Value = "More Info button through JS  not clearing old for and displaying new"
$ws.Range("C4").Value = "CSS/HTML work around "

# D4 gets the next day's date, formatted like D3 (copy D3's format then overwrite the value)
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = 43486

# Move the selection to D5 to match the saved workbook state
$ws.Range("D5").Select()
